$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 20 (2025Q2) with refreshed ADD metrics
$ws.Range("C20").Value = 337
$ws.Range("D20").Value = 264
$ws.Range("E20").Value = 73
$ws.Range("F20").Value = 81.73374613003097
